# finalizada automatizacion caso altero y caso exitoso
# Adds a new data row (row 3) to the "Datos" sheet, reproducing the
# "caso exitoso" (successful case) test data row, based on a copy of
# row 2's formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting (borders, number format, etc.) of the existing
# data row (row 2) down into the new row 3.
$ws.Range("A2:H2").Copy()
$ws.Range("A3:H3").PasteSpecial(-4122)
$ws.Rows(3).RowHeight = 19.4

# Fill in the new row's values. The order mirrors the order the cells
# were actually edited (new/changed values first, then the ones that
# reuse already-existing shared strings).
$ws.Range("B3").Value = "compan2"
$ws.Range("C3").Value = "company2"
$ws.Range("G3").Value = "error"
$ws.Range("A3").Value = "2"
$ws.Range("D3").Value = "1234567890 "
$ws.Range("E3").Value = "NENE"
$ws.Range("F3").Value = "USA"
$ws.Range("H3").Value = "1000"

# Move the active selection to the new row, as left by the editor.
$ws.Range("A3").Select() | Out-Null
